$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts Late/heading/Outstanding right)
$ws.Columns("N").Insert()

# The newly inserted column inherits the width of the column to its left (M)
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet / selection, instead of "Transactions"
$ws.Activate()
$ws.Range("S9").Select()
